$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its exact text representation
# (avoids Excel auto-converting numeric-looking strings to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '29.373.13'
$ws.Range('E2').Value = '  -0.15%  '

# Row 3
$ws.Range('D3').Value = '1.845.57'
$ws.Range('E3').Value = '  -0.26%  '

# Row 4
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').Value = '238.86'
$ws.Range('E5').Value = '  -1.40%  '

# Row 6
$ws.Range('D6').Value = '0.6293'
$ws.Range('E6').Value = '  +0.00%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').Value = '0.07554'
$ws.Range('E8').Value = '  -0.58%  '

# Row 9
$ws.Range('D9').Value = '0.2946'
$ws.Range('E9').Value = '  -0.97%  '

# Row 10
$ws.Range('D10').Value = '24.61'
$ws.Range('E10').Value = '  +0.71%  '

# Row 11
$ws.Range('D11').Value = '0.07699'
$ws.Range('E11').Value = '  -0.28%  '

# Row 12
$ws.Range('D12').Value = '1.842.89'
$ws.Range('E12').Value = '  -3.86%  '

# Row 13
$ws.Range('D13').Value = '4.980'
$ws.Range('E13').Value = '  -0.54%  '

# Row 14
$ws.Range('D14').Value = '0.6780'
$ws.Range('E14').Value = '  -1.70%  '

# Row 15
$ws.Range('D15').Value = '0.00001020'
$ws.Range('E15').Value = '  +2.43%  '

# Row 16
$ws.Range('D16').Value = '83.11'
$ws.Range('E16').Value = '  -0.18%  '

# Row 17
$ws.Range('D17').Value = '2.102.09'
$ws.Range('E17').Value = '  -4.02%  '

# Row 18
$ws.Range('D18').Value = '6.138'
$ws.Range('E18').Value = '  -0.92%  '

# Row 19
$ws.Range('D19').Value = '29.409.92'
$ws.Range('E19').Value = '  -0.40%  '

# Row 20
$ws.Range('D20').Value = '228.29'
$ws.Range('E20').Value = '  -1.97%  '

# Row 21
$ws.Range('E21').Value = '  -1.07%  '

# Row 22
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  -0.06%  '

# Row 23
$ws.Range('D23').Value = '7.461'
$ws.Range('E23').Value = '  -2.93%  '

# Row 24
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.04%  '

# Row 25
$ws.Range('D25').Value = '156.83'
$ws.Range('E25').Value = '  +1.14%  '

# Row 26
$ws.Range('D26').Value = '0.1392'
$ws.Range('E26').Value = '  -0.33%  '

# Row 27
$ws.Range('D27').Value = '8.354'
$ws.Range('E27').Value = '  -1.42%  '

# Row 28
$ws.Range('D28').Value = '17.63'
$ws.Range('E28').Value = '  -0.22%  '

# Row 29
$ws.Range('D29').Value = '1.456'
$ws.Range('E29').Value = '  -1.32%  '

# Row 30
$ws.Range('D30').Value = '1.268'
$ws.Range('E30').Value = '  +0.84%  '

# Row 31
$ws.Range('D31').Value = '0.05620'
$ws.Range('E31').Value = '  -2.78%  '

# Row 32
$ws.Range('D32').Value = '4.119'
$ws.Range('E32').Value = '  -0.27%  '

# Row 33
$ws.Range('D33').Value = '4.037'
$ws.Range('E33').Value = '  +0.38%  '

# Row 35
$ws.Range('D35').Value = '1.154'
$ws.Range('E35').Value = '  -0.81%  '

# Row 36
$ws.Range('D36').Value = '0.7139'
$ws.Range('E36').Value = '  -1.04%  '

# Row 37
$ws.Range('D37').Value = '2.593'
$ws.Range('E37').Value = '  +0.22%  '

# Row 38
$ws.Range('D38').Value = '1.241.67'
$ws.Range('E38').Value = '  -0.49%  '

# Row 39
$ws.Range('D39').Value = '0.01807'
$ws.Range('E39').Value = '  +0.15%  '

# Row 40
$ws.Range('D40').Value = '2.767'
$ws.Range('E40').Value = '  -0.95%  '

# Row 41
$ws.Range('D41').Value = '6.209'
$ws.Range('E41').Value = '  +1.93%  '

# Row 42
$ws.Range('D42').Value = '0.9029'
$ws.Range('E42').Value = '  -0.76%  '

# Row 43
$ws.Range('D43').Value = '1.0000'
$ws.Range('E43').Value = '  +0.06%  '

# Row 44
$ws.Range('D44').Value = '101.78'
$ws.Range('E44').Value = '  +0.10%  '

# Row 45
$ws.Range('D45').Value = '66.01'
$ws.Range('E45').Value = '  -2.63%  '

# Row 46
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.114'
$ws.Range('E46').Value = '  -2.48%  '

# Row 47
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.00000000118'
$ws.Range('E47').Value = '  -2.49%  '

# Row 48
$ws.Range('D48').Value = '0.3999'
$ws.Range('E48').Value = '  -0.84%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '8.991'
$ws.Range('E49').Value = '  -2.09%  '

# Row 50
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '1.683'
$ws.Range('E50').Value = '  -1.09%  '

# Row 51
$ws.Range('D51').Value = '0.1117'
$ws.Range('E51').Value = '  -0.35%  '
